# Remove responsive design implementation
#
# Appends one new trailing data row (row 95) to each of the four
# worksheets (MID_LFT_#1, MID_LFT_#2, MID_PLT_#1, MID_PLT_#2). The new
# row mirrors the existing last row (row 94) in every column except
# column A, which advances to the next day's timestamp.

$wb = $excel.ActiveWorkbook

$lastRow = 94
$newRow  = 95

# Per-sheet payload for the new row, in column order A..I.
$rowsToAdd = @(
    @{
        Sheet = 1
        A = 45881.46170138889
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,"
        D = "0x01,0x18"
        E = "0x07"
        F = 400
        G = [double]"5.68631262647113e+23"
        H = 280
        I = 7
    },
    @{
        Sheet = 2
        A = 45881.46170138889
        B = "0x01,0x7c"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
        D = "0x01,0x24"
        E = "0x19"
        F = 380
        G = [double]"5.68432987514711e+23"
        H = 292
        I = 25
    },
    @{
        Sheet = 3
        A = 45881.46170138889
        B = "0x00,0x6e"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"
        D = "0x00,0x5E"
        E = "0x15"
        F = 110
        G = [double]"5.68631262647113e+23"
        H = 94
        I = 15
    },
    @{
        Sheet = 4
        A = 45881.46170138889
        B = "0x00,0x82"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"
        D = "0x00,0x75"
        E = "0x9"
        F = 130
        G = [double]"5.68631262647113e+23"
        H = 117
        I = 9
    }
)

foreach ($entry in $rowsToAdd) {
    $ws = $wb.Worksheets.Item($entry.Sheet)

    $ws.Cells.Item($newRow, 1).Value = $entry.A
    $ws.Cells.Item($newRow, 2).Value = $entry.B
    $ws.Cells.Item($newRow, 3).Value = $entry.C
    $ws.Cells.Item($newRow, 4).Value = $entry.D
    $ws.Cells.Item($newRow, 5).Value = $entry.E
    $ws.Cells.Item($newRow, 6).Value = $entry.F
    $ws.Cells.Item($newRow, 7).Value = $entry.G
    $ws.Cells.Item($newRow, 8).Value = $entry.H
    $ws.Cells.Item($newRow, 9).Value = $entry.I

    # Column A carries the "yyyy-mm-dd h:mm:ss"-style date format used by
    # every other row in the column; copy it from the prior last row so
    # the new cell keeps the same style index instead of falling back to
    # the default General format.
    $ws.Cells.Item($newRow, 1).NumberFormat = $ws.Cells.Item($lastRow, 1).NumberFormat
}
